$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "images/westminster_abbey_a.jpg" row (row 9); everything below shifts up.
$ws.Rows.Item(9).Delete()

# After the shift, the trailing rows that used to hold
# wrigley_field_b.jpg / beetle.png / john_stamos.jpg / bird_parrot.png
# now live at rows 12-15 - remove them entirely.
$ws.Range("A12:A15").EntireRow.Delete()

# Update the active selection to match the saved state.
$ws.Range("J14").Select()
